$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The post "「全ては結びついている」" (row 517) was removed from the blog-post
# log. Deleting its entire row shifts every following row up by one, which
# is exactly what the workbook diff shows (row 518 "「初心者の..." becomes the
# new row 517, and so on through the former row 607, which becomes row 606).
$ws.Rows.Item(517).Delete()
